$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C19").NumberFormat = "@"

$ws.Range("A19").Value = "Szörényi Zalán András"
$ws.Range("B19").Value = "szorenyi.za@gmail.com"
$ws.Range("C19").Value = "10"
$ws.Range("D19").Value = "B"
$ws.Range("E19").Value = "Asztalitenisz(forgó)"
$ws.Range("F19").Value = "{}"
$ws.Range("G19").Value = "Nincs"
